# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (constant
# value "stock") right before the existing "date" column. The existing
# date / legislator_name / legislator_id columns all shift one column to
# the right (H->I, I->J, J->K). A stray space in one of the company names
# is also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a new column at H; this shifts old H:J (date, legislator_name,
# legislator_id) to I:K and widens the used range to A1:K14 automatically.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 8).Value2 = "property_category"

# Every stock record is categorized as "stock".
$lastRow = 14
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value2 = "stock"
}

# Fix the stray internal space in this company's name.
$ws.Range("B12").Value2 = "遊戲橘子數位科技股份有限公司"
